$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1689002548051
$ws.Range("C2").Value = 1689002848051
